$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C: header (t+3 factor) plus 10 data rows, mirroring columns A/B
$ws.Range("C1").Value = 2
$ws.Range("C2").Value = -5.041423880421429
$ws.Range("C3").Value = -1.155700669620174
$ws.Range("C4").Value = -0.07506705564893197
$ws.Range("C5").Value = -0.4238602485323116
$ws.Range("C6").Value = 0.01520034001876744
$ws.Range("C7").Value = 0.1043279679824023
$ws.Range("C8").Value = 0.1352696695087812
$ws.Range("C9").Value = 0.02787891322180851
$ws.Range("C10").Value = 0.02331057633078736
$ws.Range("C11").Value = 0.006819328375059187

# Mirror the header cell formatting (bold, boxed border, centered) from A1/B1
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
